$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E2").Value = 753.22077675
$wsSchedule.Range("F2").Value = 16.60539631283069
$wsSchedule.Range("E3").Value = 346.9494795
$wsSchedule.Range("F3").Value = 22.9463941468254

# --- Detailed sheet updates (Price column B) ---
$wsDetailed.Range("B21").Value = 36.06
$wsDetailed.Range("B22").Value = 45.92104
$wsDetailed.Range("B23").Value = 45.97441
$wsDetailed.Range("B24").Value = 36.06028
$wsDetailed.Range("B25").Value = 21.24005
$wsDetailed.Range("B26").Value = 0
$wsDetailed.Range("B27").Value = -16.1572
$wsDetailed.Range("B28").Value = -16.27493
$wsDetailed.Range("B29").Value = -16.47514
$wsDetailed.Range("B30").Value = -21.15844
$wsDetailed.Range("B31").Value = -21.40354
$wsDetailed.Range("B32").Value = -16.02
$wsDetailed.Range("B33").Value = -16.45302
$wsDetailed.Range("B34").Value = 5.34014
$wsDetailed.Range("B35").Value = 4.15594
$wsDetailed.Range("B36").Value = -9.5
$wsDetailed.Range("B37").Value = -10.45391
$wsDetailed.Range("B38").Value = -9.181150000000001
$wsDetailed.Range("B39").Value = -2.98349
$wsDetailed.Range("B40").Value = -3.00221
$wsDetailed.Range("B41").Value = 20.41263
$wsDetailed.Range("B42").Value = 27.59769
$wsDetailed.Range("B43").Value = 21.52393
$wsDetailed.Range("B44").Value = 21.52393
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("B47").Value = 57.09
$wsDetailed.Range("B48").Value = 57.06007

# --- Detailed sheet updates (Type column C: forecast -> historical) ---
$wsDetailed.Range("C23").Value = "historical"
$wsDetailed.Range("C24").Value = "historical"
$wsDetailed.Range("C25").Value = "historical"
$wsDetailed.Range("C26").Value = "historical"
